$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Knowledge Based Agent")
$ws.Activate()

# Row 3
$ws.Range("B3").Value = 67
$ws.Range("C3").Value = 1
$ws.Range("G3").Value = 14
$ws.Range("H3").Value = 973

# Row 4
$ws.Range("B4").Value = 108
$ws.Range("C4").Value = 1
$ws.Range("G4").Value = 21
$ws.Range("H4").Value = 951

# Row 5
$ws.Range("B5").Value = 23
$ws.Range("C5").Value = 0
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = -16

# Row 6
$ws.Range("B6").Value = 34
$ws.Range("G6").Value = 7
$ws.Range("H6").Value = 980

# Row 7
$ws.Range("B7").Value = 115
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("G7").Value = 24
$ws.Range("H7").Value = 957

# Row 8
$ws.Range("B8").Value = 112
$ws.Range("D8").Value = 2
$ws.Range("G8").Value = 21
$ws.Range("H8").Value = 967

# Row 9
$ws.Range("B9").Value = 49
$ws.Range("C9").Value = 1
$ws.Range("G9").Value = 10
$ws.Range("H9").Value = 967

# Row 10
$ws.Range("B10").Value = 63
$ws.Range("G10").Value = 12
$ws.Range("H10").Value = 974

# Row 11
$ws.Range("B11").Value = 44
$ws.Range("C11").Value = 1
$ws.Range("G11").Value = 9
$ws.Range("H11").Value = 976

# Row 12
$ws.Range("B12").Value = 33
$ws.Range("C12").Value = 0
$ws.Range("G12").Value = 7
$ws.Range("H12").Value = -24

# Row 13
$ws.Range("B13").Value = 10
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = -7

# Row 14
$ws.Range("B14").Value = 123
$ws.Range("D14").Value = 1
$ws.Range("G14").Value = 24
$ws.Range("H14").Value = 958

# Row 15
$ws.Range("B15").Value = 98
$ws.Range("C15").Value = 1
$ws.Range("G15").Value = 19
$ws.Range("H15").Value = 962

# Row 16
$ws.Range("B16").Value = 92
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -44

# Row 17
$ws.Range("B17").Value = 48
$ws.Range("H17").Value = -26

# Row 18
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0

# Row 19
$ws.Range("B19").Value = 10
$ws.Range("C19").Value = 0
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = -7

# Row 20
$ws.Range("B20").Value = 55
$ws.Range("C20").Value = 0
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -36

# Row 21
$ws.Range("B21").Value = 48
$ws.Range("D21").Value = 0
$ws.Range("G21").Value = 9
$ws.Range("H21").Value = 978

# Row 22
$ws.Range("B22").Value = 104
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("G22").Value = 21
$ws.Range("H22").Value = 967

# Update the selected cell to match the saved view state
$ws.Range("C10").Select()
